$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in rows 10 and 11 (previously blank placeholder rows) with new
# shop items "siege" and "speedup", mirroring rows 8 and 9.
$ws.Range("A10").Value = "ITEM_NAME"
$ws.Range("B10").Value = "siege"

$ws.Range("A11").Value = "ITEM_NAME"
$ws.Range("B11").Value = "speedup"

# Copy the style of A10 (blank placeholder row look) down onto the new
# blank placeholder row 12.
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null

# Move the active selection to the new last blank row, like Excel does
# after the last data entry.
$ws.Range("A12").Select() | Out-Null
